# NIT-8301339101.xlsx — "Elimina EC anteriores y se agregan nuevos, se
# modifica base de datos"
#
# The accounts-receivable table (rows 16-25, 10 worker/period rows) is
# replaced by a new table (rows 16-30, 15 worker/period rows) and two of
# the header figures (total overdue value + worker count) are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: the old table only went down to row 25; the new one goes
#    to row 30, so insert 5 fresh rows before the old last data row
#    (this also pushes the two footer/signature rows, formerly 30-31,
#    down to 35-36, and extends mergeCells/dimension accordingly).
# ---------------------------------------------------------------------
$ws.Rows("25:29").Insert()

# Give the 5 new rows the same cell formatting as the existing data rows
# (borders/fill/number formats) by copying row 24's formats across.
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J29").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Overwrite the whole worker/period table (rows 16-30) with the new
#    data set.
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047381092"
$ws.Range("D16").Value = "MAIQUI GUERRERO FLOREZ"
$ws.Range("E16").Value = "2308"
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1160000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "8647028"
$ws.Range("D17").Value = "YESITH ENRIQUE LUGO PEýA"
$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1033098595"
$ws.Range("D18").Value = "KEINER DAVID MORALES OCHOA"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 45066
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1033098595"
$ws.Range("D19").Value = "KEINER DAVID MORALES OCHOA"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "80504109"
$ws.Range("D20").Value = "TEOFILO MORALES LEON"
$ws.Range("E20").Value = "2505"
$ws.Range("F20").Value = 45066
$ws.Range("G20").Value = 1300000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "80504109"
$ws.Range("D21").Value = "TEOFILO MORALES LEON"
$ws.Range("E21").Value = "2504"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143385023"
$ws.Range("D22").Value = "DANIEL DE JESUS MARTINEZ TOVAR"
$ws.Range("E22").Value = "2303"
$ws.Range("F22").Value = 46400
$ws.Range("G22").Value = 1160000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143385023"
$ws.Range("D23").Value = "DANIEL DE JESUS MARTINEZ TOVAR"
$ws.Range("E23").Value = "2302"
$ws.Range("F23").Value = 46400
$ws.Range("G23").Value = 1160000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "51922981"
$ws.Range("D24").Value = "ALICIA VELOZA MORALES"
$ws.Range("E24").Value = "2505"
$ws.Range("F24").Value = 40214
$ws.Range("G24").Value = 1300000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "51922981"
$ws.Range("D25").Value = "ALICIA VELOZA MORALES"
$ws.Range("E25").Value = "2504"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "51922981"
$ws.Range("D26").Value = "ALICIA VELOZA MORALES"
$ws.Range("E26").Value = "2503"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "51922981"
$ws.Range("D27").Value = "ALICIA VELOZA MORALES"
$ws.Range("E27").Value = "2502"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "51922981"
$ws.Range("D28").Value = "ALICIA VELOZA MORALES"
$ws.Range("E28").Value = "2501"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "51922981"
$ws.Range("D29").Value = "ALICIA VELOZA MORALES"
$ws.Range("E29").Value = "2412"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "51922981"
$ws.Range("D30").Value = "ALICIA VELOZA MORALES"
$ws.Range("E30").Value = "2411"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

# ---------------------------------------------------------------------
# 3) Update the two summary figures in the header block.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 737546   # VALOR MORA
$ws.Range("C13").Value = 6        # Cant. Trabajadores
